$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so stale shared strings / formatting are not retained
$ws.Cells.Clear()

# ---------------------------------------------------------------------------
# Cell values (row-major order so the shared-string table is built in the
# same order as the target workbook)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Some info"
$ws.Range("B1").Value = "Site URL"
$ws.Range("C1").Value = "Some info"
$ws.Range("D1").Value = "Thematicity Index"
$ws.Range("E1").Value = "Some info"
$ws.Range("F1").Value = "Total Page"

$ws.Range("A2").Value = "data"
$ws.Range("B2").Value = "washingtonpost.com"
$ws.Range("C2").Value = "data"
$ws.Range("D2").Value = "will be filled"
$ws.Range("E2").Value = "data"
$ws.Range("F2").Value = "will be filled"

$ws.Range("A3").Value = "data"
$ws.Range("B3").Value = "www.businessinsider.com/"
$ws.Range("C3").Value = "data"
$ws.Range("D3").Value = "will be filled"
$ws.Range("E3").Value = "data"
$ws.Range("F3").Value = "will be filled"

$ws.Range("A4").Value = "data"
$ws.Range("B4").Value = "https://www.who.int/"
$ws.Range("C4").Value = "data"
$ws.Range("D4").Value = "will be filled"
$ws.Range("E4").Value = "data"
$ws.Range("F4").Value = "will be filled"

$ws.Range("A5").Value = "data"
$ws.Range("B5").Value = "macobserver.com"
$ws.Range("C5").Value = "data"
$ws.Range("D5").Value = "will be filled"
$ws.Range("E5").Value = "data"
$ws.Range("F5").Value = "will be filled"

# ---------------------------------------------------------------------------
# Styles - applied in the order needed so generated style indices line up
# with the target workbook (1=center, 2=bold header row default,
# 3=bold+center, 4=bold+center+green fill, 5=center+green fill)
# ---------------------------------------------------------------------------

# style index 1: plain centered alignment, used by the "data" label columns
# (the engine does not support comma multi-area ranges, so apply to each
# contiguous range individually)
$ws.Range("A2:A5").HorizontalAlignment = -4108
$ws.Range("C2:C5").HorizontalAlignment = -4108
$ws.Range("E2:E5").HorizontalAlignment = -4108

# style index 2: bold font applied to the whole header row (becomes the row
# default style / customFormat)
$ws.Rows(1).Font.Bold = $true

# style index 3: bold + centered, used by A1,C1,E1
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("E1").HorizontalAlignment = -4108

# style index 4: bold + centered + green fill, used by B1,D1,F1
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").Interior.Color = 11854022
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").Interior.Color = 11854022
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").Interior.Color = 11854022

# style index 5: centered + green fill (no bold), used by the data value
# columns B,D,F on rows 2-5
$ws.Range("B2:B5").HorizontalAlignment = -4108
$ws.Range("B2:B5").Interior.Color = 11854022
$ws.Range("D2:D5").HorizontalAlignment = -4108
$ws.Range("D2:D5").Interior.Color = 11854022
$ws.Range("F2:F5").HorizontalAlignment = -4108
$ws.Range("F2:F5").Interior.Color = 11854022

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 10 - 5/6
$ws.Columns("B").ColumnWidth = 25 - 5/6
$ws.Columns("C").ColumnWidth = 10 - 5/6
$ws.Columns("D").ColumnWidth = 17 - 5/6
$ws.Columns("E").ColumnWidth = 10 - 5/6
$ws.Columns("F").ColumnWidth = 12 - 5/6

# ---------------------------------------------------------------------------
# Page setup - fit to one page
# ---------------------------------------------------------------------------
$ws.PageSetup.Zoom = 100
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

# ---------------------------------------------------------------------------
# Comments
# ---------------------------------------------------------------------------
$ws.Range("B1").AddComment("Require column")
$ws.Range("D1").AddComment("Require column")
$ws.Range("F1").AddComment("Optional column" + [char]10 + "If there is no column, it will not be filled")
